$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39562, "Ryan Borges", "P&D", "Doenca", 5, 45087, 5477.6),
    @(66860, "Ana Sophia Vasconcelos", "Marketing", "Problemas pessoais", 5, 45087, 4049.09),
    @(24169, "Sr. André Borges", "Atendimento ao Cliente", "Viagem de negocios", 8, 45079, 7124.24),
    @(85556, "João Lucas da Conceição", "Recursos Humanos", "Doenca", 7, 45105, 5721.05),
    @(39264, "Maria Julia Nunes", "TI", "Doenca", 7, 45106, 4191.4),
    @(76874, "Matheus Câmara", "Atendimento ao Cliente", "Outros", 1, 45092, 8904.58),
    @(45321, "Rhavi Rios", "Juridico", "Consulta medica", 4, 45088, 5053.27),
    @(63720, "Ryan Pinto", "Atendimento ao Cliente", "Doenca", 1, 45084, 6516.76),
    @(72006, "Emanuelly Andrade", "Vendas", "Consulta medica", 8, 45095, 8372.34),
    @(56880, "Maria Flor Azevedo", "Recursos Humanos", "Consulta medica", 4, 45080, 5902.88)
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $ws.Cells.Item($row, 6).Value = $record[5]
    $ws.Cells.Item($row, 7).Value = $record[6]
    $row++
}
